$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for data rows 2-10
# from 45208 (2023-10-09) to 45212 (2023-10-13)
foreach ($row in 2..10) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -eq 45208) {
        $cell.Value2 = 45212
    }
}
